$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) values are stored as TEXT (not numbers) in this sheet.
# Force a text number-format before assignment so numeric-looking strings
# (e.g. "244.92") are kept as literal text rather than being coerced to
# a Double by Excel's usual "smart" cell-entry parsing.
$dCells = @{
    "D2" = "244.92"
    "D3" = "25.08"
    "D4" = "5.146"
    "D5" = "0.05638"
    "D6" = "6.522"
    "D7" = "2.980"
    "D8" = "0.8126"
    "D9" = "0.8377"
    "D10" = "0.1336"
    "D11" = "0.06951"
    "D12" = "0.02839"
    "D13" = "0.09403"
    "D14" = "0.001518"
    "D15" = "0.0005979"
    "D16" = "0.006222"
    "D17" = "3.499"
    "D18" = "2.106"
    "D19" = "0.3167"
    "D20" = "0.03314"
    "D21" = "0.1292"
    "D22" = "3.746"
    "D23" = "0.04688"
    "D27" = "0.00009699"
    "D28" = "0.0001899"
    "D40" = "0.03621"
    "D41" = "0.006249"
    "D43" = "0.002705"
    "D44" = "0.007721"
    "D45" = "0.00005284"
}
foreach ($addr in $dCells.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $dCells[$addr]
    # Drop back to the default "Normal" style so we do not leave a stray
    # quote-prefix / text number-format style applied to the cell (matches
    # the original workbook, where these cells carry no explicit style).
    # (Done per-cell, in the same iteration as the value write, since a
    # multi-area Range("D2,D3,...").Style assignment only affects the first
    # area in this host.)
    $cell.Style = "Normal"
}

# --- Columns B, C, E (Coin name / Link / Volume label) are plain text already
# (not numeric-looking), so a normal .Value assignment keeps them as text.
$textCells = @{
    "B10" = "WazirX"
    "C10" = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
    "E10" = "9WazirXWRX"
    "B11" = "MandalaExchangeToken"
    "C11" = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
    "E11" = "10MandalaExchangeTokenMDX"
    "B12" = "BitrueCoin"
    "C12" = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
    "E12" = "11BitrueCoinBTR"
    "B13" = "BitMartToken"
    "C13" = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
    "E13" = "12BitMartTokenBMX"
    "B14" = "BitForexToken"
    "C14" = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
    "E14" = "13BitForexTokenBF"
    "B15" = "One"
    "C15" = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
    "E15" = "14OneONEWorstin24h"
    "E27" = "26NitroExNTXBestin24h"
    "E47" = "46CoinbaseStockTokenCOIN"
}
foreach ($addr in $textCells.Keys) {
    $ws.Range($addr).Value = $textCells[$addr]
}
